$wb = $excel.ActiveWorkbook

$wsBeh = $wb.Worksheets.Item("Behandlungen")
$wsRec = $wb.Worksheets.Item("Rechnungen")

# Update "Vorname" for Alpha row (was misspelled "Alpa")
$wsRec.Range("D2").Value = "Alpha"

# Update "Stadt" values to a new unicode test string for all three rows
$wsRec.Range("G2").Value = "ÄäÜüÖößstadt"
$wsRec.Range("G3").Value = "ÄäÜüÖößstadt"
$wsRec.Range("G4").Value = "ÄäÜüÖößstadt"

# Match formatting of G3/G4 to the rest of their row (same style as G2/F3/F4)
$wsRec.Range("G3").Style = $wsRec.Range("F3").Style
$wsRec.Range("G4").Style = $wsRec.Range("F4").Style

# Reproduce the selection state left behind in each sheet
$wsRec.Range("G4").Select()
$wsBeh.Range("G4,D17").Select()
$wsBeh.Range("D17").Activate()

$wsRec.Activate()
